# ----------------------------------------------------------------------
# Product Backlog update: add Sprint 3 rows, renumber priorities, rename
# two user stories, and highlight the Sprint-0 rows with a new fill.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New story text is typed in first so it lands in the shared-string
#     table in the same order a human editing top-to-bottom would produce ---
$ws.Cells.Item(12,3).Value = 'Edit a Board for a Project'
$ws.Cells.Item(10,3).Value = 'Update Stage of my assigned Task'

# --- Rewrite the data rows (2-18) with the updated backlog content ---
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 'Admin'
$ws.Cells.Item(2,3).Value = 'View, Create, and Delete Employee'
$ws.Cells.Item(2,4).Value = 'I can manage the set of employees in the system'

$ws.Cells.Item(3,1).Value = 0
$ws.Cells.Item(3,2).Value = 'Admin'
$ws.Cells.Item(3,3).Value = 'Promote/Demote Employee to/from Manager for a Group'
$ws.Cells.Item(3,4).Value = 'I can manage who is supervising a particular group of employees'

$ws.Cells.Item(4,1).Value = 0
$ws.Cells.Item(4,2).Value = 'Manager*'
$ws.Cells.Item(4,3).Value = 'Assign/Remove Group to/from Project'
$ws.Cells.Item(4,4).Value = 'I can control who is actively collaborating to complete a project'

$ws.Cells.Item(5,1).Value = 0
$ws.Cells.Item(5,2).Value = 'Manager*'
$ws.Cells.Item(5,3).Value = 'Edit Title and/or Description of Project'
$ws.Cells.Item(5,4).Value = 'I can assign a name and description for the project to help others understand the purpose of the project'

$ws.Cells.Item(6,1).Value = 0.5
$ws.Cells.Item(6,2).Value = 'Admin'
$ws.Cells.Item(6,3).Value = 'View, Create, and Delete Groups'
$ws.Cells.Item(6,4).Value = 'I can organize employees into collaborative groups with a manager supervising work for the group'

$ws.Cells.Item(7,1).Value = 0.5
$ws.Cells.Item(7,2).Value = 'Manager*'
$ws.Cells.Item(7,3).Value = 'View, Create, and Delete Projects'
$ws.Cells.Item(7,4).Value = 'I can assign work to groups'

$ws.Cells.Item(8,1).Value = 1
$ws.Cells.Item(8,2).Value = 'Employee**'
$ws.Cells.Item(8,3).Value = 'View Board for Project'
$ws.Cells.Item(8,4).Value = 'I can overview work in progress and overall status for a project'

$ws.Cells.Item(9,1).Value = 1
$ws.Cells.Item(9,2).Value = 'Employee**'
$ws.Cells.Item(9,3).Value = 'Create a Task on the Board of a Project for my Group'
$ws.Cells.Item(9,4).Value = 'I can define work needed for a Project for my Group'

$ws.Cells.Item(10,1).Value = 1
$ws.Cells.Item(10,2).Value = 'Employee**'
$ws.Cells.Item(10,3).Value = 'Update Stage of my assigned Task'
$ws.Cells.Item(10,4).Value = 'I can notify others about the progress on a Task'

$ws.Cells.Item(11,1).Value = 2
$ws.Cells.Item(11,2).Value = 'Employee**'
$ws.Cells.Item(11,3).Value = 'Post a comment to a Task on the Board of a Project for my Group'
$ws.Cells.Item(11,4).Value = 'I can discuss the work for a Task in a space that will be preserved for later review'

$ws.Cells.Item(12,1).Value = 1
$ws.Cells.Item(12,2).Value = 'Manager*'
$ws.Cells.Item(12,3).Value = 'Edit a Board for a Project'
$ws.Cells.Item(12,4).Value = 'I can define the Stages of work expected for a project and setup for any handoffs needed during project work'

$ws.Cells.Item(13,1).Value = 1
$ws.Cells.Item(13,2).Value = 'Manager*'
$ws.Cells.Item(13,3).Value = 'Assign a Task on a Board from a Project for my Group to an Employee in my Group'
$ws.Cells.Item(13,4).Value = 'I can delegate work to an Employee in my Group'

$ws.Cells.Item(14,1).Value = 1
$ws.Cells.Item(14,2).Value = 'Manager*'
$ws.Cells.Item(14,3).Value = 'Assign/Remove Group(s) to/from Stage of Project'
$ws.Cells.Item(14,4).Value = 'I can define who will work on the project in a particular stage'

$ws.Cells.Item(15,1).Value = 3
$ws.Cells.Item(15,2).Value = 'Manager*'
$ws.Cells.Item(15,3).Value = 'View Task Status Report for Employee in my Group'
$ws.Cells.Item(15,4).Value = 'I can monitor the progress for a single Employee'

$ws.Cells.Item(16,1).Value = 3
$ws.Cells.Item(16,2).Value = 'Manager*'
$ws.Cells.Item(16,3).Value = 'View Task Status Report for a Project I Manage'
$ws.Cells.Item(16,4).Value = 'I can monitor the progress for a Project'

$ws.Cells.Item(17,1).Value = 3
$ws.Cells.Item(17,2).Value = 'Manager*'
$ws.Cells.Item(17,3).Value = 'View Task Status Report for a Group I manage'
$ws.Cells.Item(17,4).Value = 'I can monitor the progress for a single Group'

$ws.Cells.Item(18,1).Value = 4
$ws.Cells.Item(18,2).Value = 'Manager*'
$ws.Cells.Item(18,3).Value = 'Assign an Employee in my Group as a reporter for a Task on a Board for my Group'
$ws.Cells.Item(18,4).Value = 'I can define who will be notified about updates for a task'

# --- Row heights: rows that now wrap onto two lines need the taller
#     (~29pt) height; rows that no longer wrap go back to the default. ---
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(17).AutoFit()

$ws.Rows.Item(3).RowHeight = 29
$ws.Rows.Item(5).RowHeight = 29
$ws.Rows.Item(6).RowHeight = 29
$ws.Rows.Item(11).RowHeight = 29
$ws.Rows.Item(12).RowHeight = 29
$ws.Rows.Item(13).RowHeight = 29
$ws.Rows.Item(18).RowHeight = 29
$ws.Rows.Item(19).RowHeight = 29
$ws.Rows.Item(20).RowHeight = 29

# --- Highlight the new Sprint-0 rows (2-7) with a light fill + wrap ---
$ws.Range("A2:D7").Interior.ThemeColor = 4
$ws.Range("A2:D7").WrapText = $true

# --- Refresh the saved sort state to cover the extra rows (now 20) ---
$sortRange = $ws.Range("A2:D20")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A20")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("B2:B20")) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Apply()

# --- Restore the active selection to where the edit left off ---
$ws.Range("C11").Select()

